$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.357.20'
$ws.Range('E2').Value = '  +2.24%  '
$ws.Range('D3').Value = '1.648.25'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.57'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.518'
$ws.Range('E6').Value = '  +2.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.256'
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0630'
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.19'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('D12').Value = '1.879.41'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = '1.642.91'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.16'
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.548'
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.40'
$ws.Range('E16').Value = '  +1.66%  '
$ws.Range('D17').Value = '27.341.27'
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').Value = '0.0₃0743'
$ws.Range('E18').Value = '  +2.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '221.02'
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.97'
$ws.Range('E21').Value = '  +4.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.55'
$ws.Range('E22').Value = '  +4.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.43'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.21'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.33'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.51'
$ws.Range('E26').Value = '  +1.77%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.81'
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0514'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.41'
$ws.Range('E32').Value = '  +2.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.04'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.59'
$ws.Range('E34').Value = '  +2.17%  '
$ws.Range('D35').Value = '1.297.63'
$ws.Range('E35').Value = '  +3.45%  '
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.554'
$ws.Range('E38').Value = '  +4.10%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.869'
$ws.Range('E39').Value = '  +4.33%  '
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.812'
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('E42').Value = '  +6.34%  '
$ws.Range('E43').Value = '  -1.92%  '
$ws.Range('D44').Value = '1.788.90'
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.45'
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.41'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('E47').Value = '  +3.11%  '
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0514'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.76'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0975'
$ws.Range('E51').Value = '  +1.12%  '
